$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Cells.Item(2, 3).Value = 7.825809040780897
$ws.Cells.Item(2, 4).Value = 8.088210603204612
$ws.Cells.Item(2, 5).Value = 12.1568171631784
$ws.Cells.Item(2, 6).Value = 41.01905112221466
$ws.Cells.Item(2, 7).Value = 3.721019344694178
$ws.Cells.Item(2, 9).Value = 35.53810612300434
$ws.Cells.Item(2, 11).Value = 22.30392397076362
$ws.Cells.Item(2, 12).Value = 9.860337436749228
$ws.Cells.Item(2, 14).Value = 19.38035546254005
$ws.Cells.Item(3, 3).Value = 7.786805768057308
$ws.Cells.Item(3, 4).Value = 8.126481819646244
$ws.Cells.Item(3, 5).Value = 12.1222666669897
$ws.Cells.Item(3, 6).Value = 40.89232976696306
$ws.Cells.Item(3, 7).Value = 3.725220503346018
$ws.Cells.Item(3, 9).Value = 35.49398967238753
$ws.Cells.Item(3, 11).Value = 21.91405395893464
$ws.Cells.Item(3, 12).Value = 9.867151596782078
$ws.Cells.Item(3, 14).Value = 19.46180322354031
$ws.Cells.Item(4, 3).Value = 7.764255731764505
$ws.Cells.Item(4, 4).Value = 8.151068502986792
$ws.Cells.Item(4, 5).Value = 12.10377882119274
$ws.Cells.Item(4, 6).Value = 40.82738261361273
$ws.Cells.Item(4, 7).Value = 3.727931403224259
$ws.Cells.Item(4, 9).Value = 35.47711145676281
$ws.Cells.Item(4, 11).Value = 21.6773292543424
$ws.Cells.Item(4, 12).Value = 9.873195927503135
$ws.Cells.Item(4, 14).Value = 19.51387244045306
$ws.Cells.Item(5, 3).Value = 7.755424562167018
$ws.Cells.Item(5, 4).Value = 8.161362066386097
$ws.Cells.Item(5, 5).Value = 12.09693438141958
$ws.Cells.Item(5, 6).Value = 40.80415665706688
$ws.Cells.Item(5, 7).Value = 3.72906928590698
$ws.Cells.Item(5, 9).Value = 35.47279665164541
$ws.Cells.Item(5, 11).Value = 21.58166150025873
$ws.Cells.Item(5, 12).Value = 9.876126210741679
$ws.Cells.Item(5, 14).Value = 19.53561127435258
$ws.Cells.Item(6, 3).Value = 7.75397996783849
$ws.Cells.Item(6, 4).Value = 8.163087892036094
$ws.Cells.Item(6, 5).Value = 12.09583961356272
$ws.Cells.Item(6, 6).Value = 40.80049583641698
$ws.Cells.Item(6, 7).Value = 3.729260237729994
$ws.Cells.Item(6, 9).Value = 35.47223482896761
$ws.Cells.Item(6, 11).Value = 21.56582806080312
$ws.Cells.Item(6, 12).Value = 9.876640976234635
$ws.Cells.Item(6, 14).Value = 19.53925247500354
$ws.Cells.Item(7, 3).Value = 7.764135173014462
$ws.Cells.Item(7, 4).Value = 8.151206214103574
$ws.Cells.Item(7, 5).Value = 12.10368371819107
$ws.Cells.Item(7, 6).Value = 40.82705625409763
$ws.Cells.Item(7, 7).Value = 3.72794661463208
$ws.Cells.Item(7, 9).Value = 35.47704289479759
$ws.Cells.Item(7, 11).Value = 21.6760356352806
$ws.Cells.Item(7, 12).Value = 9.873233555745044
$ws.Cells.Item(7, 14).Value = 19.51416350884221
$ws.Cells.Item(8, 3).Value = 7.812074617957689
$ws.Cells.Item(8, 4).Value = 8.101181127285399
$ws.Cells.Item(8, 5).Value = 12.14434089363224
$ws.Cells.Item(8, 6).Value = 40.97269376764717
$ws.Cells.Item(8, 7).Value = 3.722440719244072
$ws.Cells.Item(8, 9).Value = 35.52077404576572
$ws.Cells.Item(8, 11).Value = 22.16902412506573
$ws.Cells.Item(8, 12).Value = 9.862300482703516
$ws.Cells.Item(8, 14).Value = 19.40801244498535
$ws.Cells.Item(9, 3).Value = 7.916867723206636
$ws.Cells.Item(9, 4).Value = 8.011682434157013
$ws.Cells.Item(9, 5).Value = 12.24549696811344
$ws.Cells.Item(9, 6).Value = 41.35983907863763
$ws.Cells.Item(9, 7).Value = 3.71267990587754
$ws.Cells.Item(9, 9).Value = 35.68761848805428
$ws.Cells.Item(9, 11).Value = 23.15110717954797
$ws.Cells.Item(9, 12).Value = 9.855646466389775
$ws.Cells.Item(9, 14).Value = 19.21609247418676
$ws.Cells.Item(10, 3).Value = 7.99999144788069
$ws.Cells.Item(10, 4).Value = 7.95112488191588
$ws.Cells.Item(10, 5).Value = 12.33256185038842
$ws.Cells.Item(10, 6).Value = 41.70524897749265
$ws.Cells.Item(10, 7).Value = 3.706131743859007
$ws.Cells.Item(10, 9).Value = 35.85959709256677
$ws.Cells.Item(10, 11).Value = 23.87441591740145
$ws.Cells.Item(10, 12).Value = 9.859797741105657
$ws.Cells.Item(10, 14).Value = 19.08484334903934
$ws.Cells.Item(11, 3).Value = 8.03903613631446
$ws.Cells.Item(11, 4).Value = 7.924694634363116
$ws.Cells.Item(11, 5).Value = 12.37485587552046
$ws.Cells.Item(11, 6).Value = 41.87533854952519
$ws.Cells.Item(11, 7).Value = 3.703286273797907
$ws.Cells.Item(11, 9).Value = 35.9485033739465
$ws.Cells.Item(11, 11).Value = 24.20230923975707
$ws.Cells.Item(11, 12).Value = 9.863651717950221
$ws.Cells.Item(11, 14).Value = 19.02722108241564
$ws.Cells.Item(12, 3).Value = 8.053989285127974
$ws.Cells.Item(12, 4).Value = 7.914846239020486
$ws.Cells.Item(12, 5).Value = 12.3912500693438
$ws.Cells.Item(12, 6).Value = 41.9415802492818
$ws.Cells.Item(12, 7).Value = 3.702227797444926
$ws.Cells.Item(12, 9).Value = 35.98369554813085
$ws.Cells.Item(12, 11).Value = 24.32618676748814
$ws.Cells.Item(12, 12).Value = 9.865393642415402
$ws.Cells.Item(12, 14).Value = 19.00569824787824
$ws.Cells.Item(13, 3).Value = 8.050761545945871
$ws.Cells.Item(13, 4).Value = 7.916960152376451
$ws.Cells.Item(13, 5).Value = 12.38770259443329
$ws.Cells.Item(13, 6).Value = 41.92723295263671
$ws.Cells.Item(13, 7).Value = 3.702454914462898
$ws.Cells.Item(13, 9).Value = 35.97604862058277
$ws.Cells.Item(13, 11).Value = 24.29952227717837
$ws.Cells.Item(13, 12).Value = 9.865005927187793
$ws.Cells.Item(13, 14).Value = 19.01032037972949
$ws.Cells.Item(14, 3).Value = 8.04026303259149
$ws.Cells.Item(14, 4).Value = 7.923881194107113
$ws.Cells.Item(14, 5).Value = 12.37619710145741
$ws.Cells.Item(14, 6).Value = 41.88075176125765
$ws.Cells.Item(14, 7).Value = 3.703198811437787
$ws.Cells.Item(14, 9).Value = 35.95136814835904
$ws.Cells.Item(14, 11).Value = 24.21250719320055
$ws.Cells.Item(14, 12).Value = 9.86378936663675
$ws.Cells.Item(14, 14).Value = 19.02544443652285
$ws.Cells.Item(15, 3).Value = 8.033853958787599
$ws.Cells.Item(15, 4).Value = 7.928141374152609
$ws.Cells.Item(15, 5).Value = 12.36919867729384
$ws.Cells.Item(15, 6).Value = 41.85251833785198
$ws.Cells.Item(15, 7).Value = 3.703656945709903
$ws.Cells.Item(15, 9).Value = 35.9364489656239
$ws.Cells.Item(15, 11).Value = 24.15916677200861
$ws.Cells.Item(15, 12).Value = 9.863080972035998
$ws.Cells.Item(15, 14).Value = 19.03474703307662
$ws.Cells.Item(16, 3).Value = 7.997463774963371
$ws.Cells.Item(16, 4).Value = 7.952874603457134
$ws.Cells.Item(16, 5).Value = 12.32985123985016
$ws.Cells.Item(16, 6).Value = 41.694391442564
$ws.Cells.Item(16, 7).Value = 3.706320373665994
$ws.Cells.Item(16, 9).Value = 35.85400086648677
$ws.Cells.Item(16, 11).Value = 23.85295357497521
$ws.Cells.Item(16, 12).Value = 9.859585418648331
$ws.Cells.Item(16, 14).Value = 19.08865084321218
$ws.Cells.Item(17, 3).Value = 7.975448421862105
$ws.Cells.Item(17, 4).Value = 7.968333502738756
$ws.Cells.Item(17, 5).Value = 12.30639567228441
$ws.Cells.Item(17, 6).Value = 41.6006842649797
$ws.Cells.Item(17, 7).Value = 3.707988354679136
$ws.Cells.Item(17, 9).Value = 35.80614916883372
$ws.Cells.Item(17, 11).Value = 23.66471883135646
$ws.Cells.Item(17, 12).Value = 9.857944389270944
$ws.Cells.Item(17, 14).Value = 19.12225121614604
$ws.Cells.Item(18, 3).Value = 7.962902242548796
$ws.Cells.Item(18, 4).Value = 7.977330286748989
$ws.Cells.Item(18, 5).Value = 12.29315816300386
$ws.Cells.Item(18, 6).Value = 41.54800843733132
$ws.Cells.Item(18, 7).Value = 3.708960289872815
$ws.Cells.Item(18, 9).Value = 35.77963138658726
$ws.Cells.Item(18, 11).Value = 23.55635066472336
$ws.Cells.Item(18, 12).Value = 9.857185537443168
$ws.Cells.Item(18, 14).Value = 19.14177351083783
$ws.Cells.Item(19, 3).Value = 7.958674603880479
$ws.Cells.Item(19, 4).Value = 7.980394534517727
$ws.Cells.Item(19, 5).Value = 12.28871994868321
$ws.Cells.Item(19, 6).Value = 41.53038408788556
$ws.Cells.Item(19, 7).Value = 3.709291531008578
$ws.Cells.Item(19, 9).Value = 35.77082581957599
$ws.Cells.Item(19, 11).Value = 23.51964567052612
$ws.Cells.Item(19, 12).Value = 9.856960382799025
$ws.Cells.Item(19, 14).Value = 19.14841719592194
$ws.Cells.Item(20, 3).Value = 7.977780006017428
$ws.Cells.Item(20, 4).Value = 7.966676989297857
$ws.Cells.Item(20, 5).Value = 12.30886637687532
$ws.Cells.Item(20, 6).Value = 41.61053328975257
$ws.Cells.Item(20, 7).Value = 3.707809496632789
$ws.Cells.Item(20, 9).Value = 35.81113908462673
$ws.Cells.Item(20, 11).Value = 23.68476806099112
$ws.Cells.Item(20, 12).Value = 9.858099930728006
$ws.Cells.Item(20, 14).Value = 19.11865410369245
$ws.Cells.Item(21, 3).Value = 8.043342222481773
$ws.Cells.Item(21, 4).Value = 7.921843973857285
$ws.Cells.Item(21, 5).Value = 12.37956634275287
$ws.Cells.Item(21, 6).Value = 41.89435494923593
$ws.Cells.Item(21, 7).Value = 3.70297979512074
$ws.Cells.Item(21, 9).Value = 35.95857609185144
$ws.Cells.Item(21, 11).Value = 24.23807442613748
$ws.Cells.Item(21, 12).Value = 9.864139034714773
$ws.Cells.Item(21, 14).Value = 19.02099407841655
$ws.Cells.Item(22, 3).Value = 8.087164341118832
$ws.Cells.Item(22, 4).Value = 7.893476331808785
$ws.Cells.Item(22, 5).Value = 12.42797406502755
$ws.Cells.Item(22, 6).Value = 42.09051271931508
$ws.Cells.Item(22, 7).Value = 3.69993423367122
$ws.Cells.Item(22, 9).Value = 36.06381925536685
$ws.Cells.Item(22, 11).Value = 24.59795657999098
$ws.Cells.Item(22, 12).Value = 9.86973242367727
$ws.Cells.Item(22, 14).Value = 18.9589004736433
$ws.Cells.Item(23, 3).Value = 8.063689718540529
$ws.Cells.Item(23, 4).Value = 7.908531455350873
$ws.Cells.Item(23, 5).Value = 12.40193938290959
$ws.Cells.Item(23, 6).Value = 41.98485517866533
$ws.Cells.Item(23, 7).Value = 3.701549600157884
$ws.Cells.Item(23, 9).Value = 36.00683977510511
$ws.Cells.Item(23, 11).Value = 24.40607907489156
$ws.Cells.Item(23, 12).Value = 9.866596561512244
$ws.Cells.Item(23, 14).Value = 18.99188315331411
$ws.Cells.Item(24, 3).Value = 7.976725551188321
$ws.Cells.Item(24, 4).Value = 7.967425558671392
$ws.Cells.Item(24, 5).Value = 12.30774860009651
$ws.Cells.Item(24, 6).Value = 41.6060768128294
$ws.Cells.Item(24, 7).Value = 3.707890317882566
$ws.Cells.Item(24, 9).Value = 35.80888005067627
$ws.Cells.Item(24, 11).Value = 23.67570426047152
$ws.Cells.Item(24, 12).Value = 9.858029035445849
$ws.Cells.Item(24, 14).Value = 19.12027971956335
$ws.Cells.Item(25, 3).Value = 7.887407951374116
$ws.Cells.Item(25, 4).Value = 8.034978282091757
$ws.Cells.Item(25, 5).Value = 12.21586761562084
$ws.Cells.Item(25, 6).Value = 41.24431387743122
$ws.Cells.Item(25, 7).Value = 3.715210419715592
$ws.Cells.Item(25, 9).Value = 35.63379566105591
$ws.Cells.Item(25, 11).Value = 22.88457887860968
$ws.Cells.Item(25, 12).Value = 9.855860309069987
$ws.Cells.Item(25, 14).Value = 19.26628830543068
